$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 916.72
$ws.Range("I33").Value = 732.94116
$ws.Range("J33").Value = 1307.25
$ws.Range("K33").Value = 732.94116
$ws.Range("L33").Value = 1307.25
$ws.Range("M33").Value = -503.94116
$ws.Range("N33").Value = -1765.25
$ws.Range("H55").Value = 125001420
$ws.Range("I55").Value = 250002540
$ws.Range("K55").Value = 250002540
$ws.Range("M55").Value = -250002326
$ws.Range("H70").Value = 2842088
$ws.Range("I70").Value = 3247886.5
$ws.Range("J70").Value = 1500
$ws.Range("K70").Value = 9743659.5
$ws.Range("L70").Value = 4500
$ws.Range("M70").Value = -9743389.5
$ws.Range("N70").Value = -5040
$ws.Range("H73").Value = 2842088
$ws.Range("I73").Value = 3247886.5
$ws.Range("J73").Value = 1500
$ws.Range("K73").Value = 9743659.5
$ws.Range("L73").Value = 4500
$ws.Range("M73").Value = -9742723.5
$ws.Range("N73").Value = -6372
$ws.Range("H98").Value = 1542.7142
$ws.Range("I98").Value = 1324.75
$ws.Range("J98").Value = 1833.3334
$ws.Range("K98").Value = 1324.75
$ws.Range("L98").Value = 1833.3334
$ws.Range("M98").Value = 173.25
$ws.Range("N98").Value = -4829.3334
$ws.Range("H107").Value = 1125.8636
$ws.Range("I107").Value = 1125.9445
$ws.Range("J107").Value = 1125.5
$ws.Range("K107").Value = 1125.9445
$ws.Range("L107").Value = 1125.5
$ws.Range("M107").Value = 794.0554999999999
$ws.Range("N107").Value = -4965.5
$ws.Range("H112").Value = 8149.8057
$ws.Range("J112").Value = 9002.906000000001
$ws.Range("L112").Value = 27008.718
$ws.Range("N112").Value = -29224.718
$ws.Range("H122").Value = 1542.7142
$ws.Range("I122").Value = 1324.75
$ws.Range("J122").Value = 1833.3334
$ws.Range("K122").Value = 3974.25
$ws.Range("L122").Value = 5500.0002
$ws.Range("M122").Value = -1524.25
$ws.Range("N122").Value = -10400.0002
$ws.Range("H129").Value = 1196.95
$ws.Range("I129").Value = 681.125
$ws.Range("J129").Value = 1325.9062
$ws.Range("K129").Value = 2043.375
$ws.Range("L129").Value = 3977.7186
$ws.Range("M129").Value = 2956.625
$ws.Range("N129").Value = -13977.7186
$ws.Range("H132").Value = 6858.7
$ws.Range("I132").Value = 6509.9443
$ws.Range("J132").Value = 9997.5
$ws.Range("K132").Value = 19529.8329
$ws.Range("L132").Value = 29992.5
$ws.Range("M132").Value = -16999.8329
$ws.Range("N132").Value = -35052.5
$ws.Range("H135").Value = 3121.8333
$ws.Range("I135").Value = 3121.8333
$ws.Range("K135").Value = 28096.4997
$ws.Range("M135").Value = -25561.4997
$ws.Range("H138").Value = 196159.2
$ws.Range("I138").Value = 2384.389
$ws.Range("J138").Value = 287947.25
$ws.Range("K138").Value = 7153.167
$ws.Range("L138").Value = 863841.75
$ws.Range("M138").Value = -2013.167
$ws.Range("N138").Value = -874121.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 543665.4
$ws.Range("I32").Value = 729453.5600000001
$ws.Range("J32").Value = 11072.533
$ws.Range("K32").Value = 729453.5600000001
$ws.Range("L32").Value = 11072.533
$ws.Range("M32").Value = -729166.5600000001
$ws.Range("N32").Value = -11646.533
$ws.Range("H80").Value = 24121.334
$ws.Range("J80").Value = 21249
$ws.Range("L80").Value = 21249
$ws.Range("N80").Value = -23245
$ws.Range("H83").Value = 24121.334
$ws.Range("J83").Value = 21249
$ws.Range("L83").Value = 63747
$ws.Range("N83").Value = -73731
$ws.Range("H102").Value = 1928.4286
$ws.Range("I102").Value = 2020
$ws.Range("J102").Value = 1699.5
$ws.Range("K102").Value = 2020
$ws.Range("L102").Value = 1699.5
$ws.Range("M102").Value = -398
$ws.Range("N102").Value = -4943.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 686.9
$ws.Range("I64").Value = 683.3333
$ws.Range("J64").Value = 692.25
$ws.Range("K64").Value = 683.3333
$ws.Range("L64").Value = 692.25
$ws.Range("M64").Value = -458.3333
$ws.Range("N64").Value = -1142.25
$ws.Range("H67").Value = 686.9
$ws.Range("I67").Value = 683.3333
$ws.Range("J67").Value = 692.25
$ws.Range("K67").Value = 683.3333
$ws.Range("L67").Value = 692.25
$ws.Range("M67").Value = 96.66669999999999
$ws.Range("N67").Value = -2252.25
$ws.Range("H86").Value = 66668468
$ws.Range("I86").Value = 83335090
$ws.Range("K86").Value = 83335090
$ws.Range("M86").Value = -83333967
$ws.Range("H89").Value = 66668468
$ws.Range("I89").Value = 83335090
$ws.Range("K89").Value = 416675450
$ws.Range("M89").Value = -416669834
$ws.Range("H107").Value = 1140.9166
$ws.Range("I107").Value = 552.1429000000001
$ws.Range("J107").Value = 1965.2
$ws.Range("K107").Value = 552.1429000000001
$ws.Range("L107").Value = 1965.2
$ws.Range("M107").Value = 1367.8571
$ws.Range("N107").Value = -5805.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1263.5454
$ws.Range("I58").Value = 829
$ws.Range("J58").Value = 1426.5
$ws.Range("K58").Value = 829
$ws.Range("L58").Value = 1426.5
$ws.Range("M58").Value = -626
$ws.Range("N58").Value = -1832.5
$ws.Range("H59").Value = 32499
$ws.Range("J59").Value = 32499
$ws.Range("L59").Value = 32499
$ws.Range("N59").Value = -34789
$ws.Range("H136").Value = 1263.5454
$ws.Range("I136").Value = 829
$ws.Range("J136").Value = 1426.5
$ws.Range("K136").Value = 2487
$ws.Range("L136").Value = 4279.5
$ws.Range("M136").Value = 63
$ws.Range("N136").Value = -9379.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 12830749
$ws.Range("I137").Value = 16678880
$ws.Range("J137").Value = 3644.3333
$ws.Range("K137").Value = 50036640
$ws.Range("L137").Value = 10932.9999
$ws.Range("M137").Value = -50031540
$ws.Range("N137").Value = -21132.9999
$ws.Range("H138").Value = 3410.5557
$ws.Range("I138").Value = 797.5
$ws.Range("J138").Value = 3737.1875
$ws.Range("K138").Value = 2392.5
$ws.Range("L138").Value = 11211.5625
$ws.Range("M138").Value = 2747.5
$ws.Range("N138").Value = -21491.5625
$ws.Range("H140").Value = 2372.8333
$ws.Range("I140").Value = 1662.6666
$ws.Range("J140").Value = 3083
$ws.Range("K140").Value = 4987.9998
$ws.Range("L140").Value = 9249
$ws.Range("M140").Value = 192.0002000000004
$ws.Range("N140").Value = -19609

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 675198.5
$ws.Range("I80").Value = 1504462.4
$ws.Range("J80").Value = 53250.625
$ws.Range("K80").Value = 1504462.4
$ws.Range("L80").Value = 53250.625
$ws.Range("M80").Value = -1503464.4
$ws.Range("N80").Value = -55246.625
$ws.Range("H83").Value = 675198.5
$ws.Range("I83").Value = 1504462.4
$ws.Range("J83").Value = 53250.625
$ws.Range("K83").Value = 7522312
$ws.Range("L83").Value = 266253.125
$ws.Range("M83").Value = -7517320
$ws.Range("N83").Value = -276237.125
$ws.Range("H132").Value = 2076.8
$ws.Range("I132").Value = 1576.909
$ws.Range("K132").Value = 4730.727000000001
$ws.Range("M132").Value = -2200.727000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 4649.276
$ws.Range("J22").Value = 8250.6
$ws.Range("L22").Value = 8250.6
$ws.Range("N22").Value = -8840.6
$ws.Range("H27").Value = 4649.276
$ws.Range("J27").Value = 8250.6
$ws.Range("L27").Value = 8250.6
$ws.Range("N27").Value = -8464.6
$ws.Range("H46").Value = 1566.6666
$ws.Range("I46").Value = 1700
$ws.Range("J46").Value = 1300
$ws.Range("K46").Value = 1700
$ws.Range("L46").Value = 1300
$ws.Range("M46").Value = -1512
$ws.Range("N46").Value = -1676
$ws.Range("H82").Value = 1798.9412
$ws.Range("I82").Value = 1724
$ws.Range("J82").Value = 1865.5555
$ws.Range("K82").Value = 1724
$ws.Range("L82").Value = 1865.5555
$ws.Range("M82").Value = -1363
$ws.Range("N82").Value = -2587.5555
$ws.Range("H85").Value = 1798.9412
$ws.Range("I85").Value = 1724
$ws.Range("J85").Value = 1865.5555
$ws.Range("K85").Value = 1724
$ws.Range("L85").Value = 1865.5555
$ws.Range("M85").Value = -476
$ws.Range("N85").Value = -4361.5555

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H92").Value = 31850
$ws.Range("J92").Value = 31850
$ws.Range("L92").Value = 31850
$ws.Range("N92").Value = -36842
$ws.Range("H132").Value = 3970753.5
$ws.Range("I132").Value = 2695.1304
$ws.Range("J132").Value = 8774193
$ws.Range("K132").Value = 8085.3912
$ws.Range("L132").Value = 26322579
$ws.Range("M132").Value = -5555.3912
$ws.Range("N132").Value = -26327639
